$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"25.35940266666667"
$ws.Range("H2").Value = [double]"76.078208"
$ws.Range("I2").Value = [double]"0.005186643687654987"
$ws.Range("J2").Value = [double]"0.005186643687654986"
$ws.Range("M2").Value = [double]"509.3923236666667"
$ws.Range("N2").Value = [double]"1528.176971"
$ws.Range("O2").Value = [double]"0.831019558191033"
$ws.Range("P2").Value = [double]"0.8310195581910331"
$ws.Range("Q2").Value = [double]"12917.885051172"
$ws.Range("R2").Value = [double]"116260.965460548"
$ws.Range("S2").Value = [double]"0.004310202345809358"
$ws.Range("T2").Value = [double]"0.004310202345809358"
$ws.Range("G3").Value = [double]"25.35940266666667"
$ws.Range("H3").Value = [double]"76.078208"
$ws.Range("I3").Value = [double]"0.005186643687654987"
$ws.Range("J3").Value = [double]"0.005186643687654986"
$ws.Range("N3").Value = [double]"3.779073"
$ws.Range("O3").Value = [double]"0.002055052284145212"
$ws.Range("P3").Value = [double]"0.002055052284145212"
$ws.Range("Q3").Value = [double]"31.945011304576"
$ws.Range("R3").Value = [double]"287.505101741184"
$ws.Range("S3").Value = [double]"1.065882395736273E-05"
$ws.Range("T3").Value = [double]"1.065882395736273E-05"
$ws.Range("G4").Value = [double]"25.35940266666667"
$ws.Range("H4").Value = [double]"76.078208"
$ws.Range("I4").Value = [double]"0.005186643687654987"
$ws.Range("J4").Value = [double]"0.005186643687654986"
$ws.Range("M4").Value = [double]"22.92703233333333"
$ws.Range("N4").Value = [double]"68.781097"
$ws.Range("O4").Value = [double]"0.03740302198339736"
$ws.Range("P4").Value = [double]"0.03740302198339737"
$ws.Range("Q4").Value = [double]"581.4158448926862"
$ws.Range("R4").Value = [double]"5232.742604034176"
$ws.Range("S4").Value = [double]"0.0001939961478694086"
$ws.Range("T4").Value = [double]"0.0001939961478694086"
$ws.Range("G5").Value = [double]"25.35940266666667"
$ws.Range("H5").Value = [double]"76.078208"
$ws.Range("I5").Value = [double]"0.005186643687654987"
$ws.Range("J5").Value = [double]"0.005186643687654986"
$ws.Range("M5").Value = [double]"79.39367866666667"
$ws.Range("N5").Value = [double]"238.181036"
$ws.Range("O5").Value = [double]"0.1295223675414243"
$ws.Range("P5").Value = [double]"0.1295223675414243"
$ws.Range("Q5").Value = [double]"2013.376266495943"
$ws.Range("R5").Value = [double]"18120.38639846349"
$ws.Range("S5").Value = [double]"0.0006717863700188576"
$ws.Range("T5").Value = [double]"0.0006717863700188575"
$ws.Range("I6").Value = [double]"0.9837462940761621"
$ws.Range("J6").Value = [double]"0.983746294076162"
$ws.Range("M6").Value = [double]"509.3923236666667"
$ws.Range("N6").Value = [double]"1528.176971"
$ws.Range("O6").Value = [double]"0.831019558191033"
$ws.Range("P6").Value = [double]"0.8310195581910331"
$ws.Range("Q6").Value = [double]"2450124.263719739"
$ws.Range("R6").Value = [double]"22051118.37347765"
$ws.Range("S6").Value = [double]"0.8175124106752383"
$ws.Range("T6").Value = [double]"0.8175124106752383"
$ws.Range("I7").Value = [double]"0.9837462940761621"
$ws.Range("J7").Value = [double]"0.983746294076162"
$ws.Range("N7").Value = [double]"3.779073"
$ws.Range("O7").Value = [double]"0.002055052284145212"
$ws.Range("P7").Value = [double]"0.002055052284145212"
$ws.Range("S7").Value = [double]"0.002021650068660604"
$ws.Range("T7").Value = [double]"0.002021650068660605"
$ws.Range("I8").Value = [double]"0.9837462940761621"
$ws.Range("J8").Value = [double]"0.983746294076162"
$ws.Range("M8").Value = [double]"22.92703233333333"
$ws.Range("N8").Value = [double]"68.781097"
$ws.Range("O8").Value = [double]"0.03740302198339736"
$ws.Range("P8").Value = [double]"0.03740302198339737"
$ws.Range("Q8").Value = [double]"110276.6484791904"
$ws.Range("R8").Value = [double]"992489.8363127136"
$ws.Range("S8").Value = [double]"0.03679508426341637"
$ws.Range("T8").Value = [double]"0.03679508426341638"
$ws.Range("I9").Value = [double]"0.9837462940761621"
$ws.Range("J9").Value = [double]"0.983746294076162"
$ws.Range("M9").Value = [double]"79.39367866666667"
$ws.Range("N9").Value = [double]"238.181036"
$ws.Range("O9").Value = [double]"0.1295223675414243"
$ws.Range("P9").Value = [double]"0.1295223675414243"
$ws.Range("Q9").Value = [double]"381875.3629559208"
$ws.Range("R9").Value = [double]"3436878.266603287"
$ws.Range("S9").Value = [double]"0.1274171490688468"
$ws.Range("T9").Value = [double]"0.1274171490688467"
$ws.Range("G10").Value = [double]"51.27300266666666"
$ws.Range("H10").Value = [double]"153.819008"
$ws.Range("I10").Value = [double]"0.01048663484403512"
$ws.Range("J10").Value = [double]"0.01048663484403512"
$ws.Range("M10").Value = [double]"509.3923236666667"
$ws.Range("N10").Value = [double]"1528.176971"
$ws.Range("O10").Value = [double]"0.831019558191033"
$ws.Range("P10").Value = [double]"0.8310195581910331"
$ws.Range("Q10").Value = [double]"26118.07396974053"
$ws.Range("R10").Value = [double]"235062.6657276648"
$ws.Range("S10").Value = [double]"0.008714598655000762"
$ws.Range("T10").Value = [double]"0.008714598655000762"
$ws.Range("G11").Value = [double]"51.27300266666666"
$ws.Range("H11").Value = [double]"153.819008"
$ws.Range("I11").Value = [double]"0.01048663484403512"
$ws.Range("J11").Value = [double]"0.01048663484403512"
$ws.Range("N11").Value = [double]"3.779073"
$ws.Range("O11").Value = [double]"0.002055052284145212"
$ws.Range("P11").Value = [double]"0.002055052284145212"
$ws.Range("Q11").Value = [double]"64.58814000217599"
$ws.Range("R11").Value = [double]"581.293260019584"
$ws.Range("S11").Value = [double]"2.155058288923115E-05"
$ws.Range("T11").Value = [double]"2.155058288923116E-05"
$ws.Range("G12").Value = [double]"51.27300266666666"
$ws.Range("H12").Value = [double]"153.819008"
$ws.Range("I12").Value = [double]"0.01048663484403512"
$ws.Range("J12").Value = [double]"0.01048663484403512"
$ws.Range("M12").Value = [double]"22.92703233333333"
$ws.Range("N12").Value = [double]"68.781097"
$ws.Range("O12").Value = [double]"0.03740302198339736"
$ws.Range("P12").Value = [double]"0.03740302198339737"
$ws.Range("Q12").Value = [double]"1175.537789965753"
$ws.Range("R12").Value = [double]"10579.84010969178"
$ws.Range("S12").Value = [double]"0.0003922318336033065"
$ws.Range("T12").Value = [double]"0.0003922318336033065"
$ws.Range("G13").Value = [double]"51.27300266666666"
$ws.Range("H13").Value = [double]"153.819008"
$ws.Range("I13").Value = [double]"0.01048663484403512"
$ws.Range("J13").Value = [double]"0.01048663484403512"
$ws.Range("M13").Value = [double]"79.39367866666667"
$ws.Range("N13").Value = [double]"238.181036"
$ws.Range("O13").Value = [double]"0.1295223675414243"
$ws.Range("P13").Value = [double]"0.1295223675414243"
$ws.Range("Q13").Value = [double]"4070.752297992477"
$ws.Range("R13").Value = [double]"36636.77068193228"
$ws.Range("S13").Value = [double]"0.001358253772541824"
$ws.Range("T13").Value = [double]"0.001358253772541824"
$ws.Range("G14").Value = [double]"2.837922333333333"
$ws.Range("H14").Value = [double]"8.513767"
$ws.Range("I14").Value = [double]"0.0005804273921477663"
$ws.Range("J14").Value = [double]"0.0005804273921477662"
$ws.Range("M14").Value = [double]"509.3923236666667"
$ws.Range("N14").Value = [double]"1528.176971"
$ws.Range("O14").Value = [double]"0.831019558191033"
$ws.Range("P14").Value = [double]"0.8310195581910331"
$ws.Range("Q14").Value = [double]"1445.615851762195"
$ws.Range("R14").Value = [double]"13010.54266585976"
$ws.Range("S14").Value = [double]"0.0004823465149846102"
$ws.Range("T14").Value = [double]"0.0004823465149846102"
$ws.Range("G15").Value = [double]"2.837922333333333"
$ws.Range("H15").Value = [double]"8.513767"
$ws.Range("I15").Value = [double]"0.0005804273921477663"
$ws.Range("J15").Value = [double]"0.0005804273921477662"
$ws.Range("N15").Value = [double]"3.779073"
$ws.Range("O15").Value = [double]"0.002055052284145212"
$ws.Range("P15").Value = [double]"0.002055052284145212"
$ws.Range("Q15").Value = [double]"3.574905221998999"
$ws.Range("R15").Value = [double]"32.174146997991"
$ws.Range("S15").Value = [double]"1.192808638013716E-06"
$ws.Range("T15").Value = [double]"1.192808638013716E-06"
$ws.Range("G16").Value = [double]"2.837922333333333"
$ws.Range("H16").Value = [double]"8.513767"
$ws.Range("I16").Value = [double]"0.0005804273921477663"
$ws.Range("J16").Value = [double]"0.0005804273921477662"
$ws.Range("M16").Value = [double]"22.92703233333333"
$ws.Range("N16").Value = [double]"68.781097"
$ws.Range("O16").Value = [double]"0.03740302198339736"
$ws.Range("P16").Value = [double]"0.03740302198339737"
$ws.Range("Q16").Value = [double]"65.06513709582211"
$ws.Range("R16").Value = [double]"585.586233862399"
$ws.Range("S16").Value = [double]"2.17097385082689E-05"
$ws.Range("T16").Value = [double]"2.17097385082689E-05"
$ws.Range("G17").Value = [double]"2.837922333333333"
$ws.Range("H17").Value = [double]"8.513767"
$ws.Range("I17").Value = [double]"0.0005804273921477663"
$ws.Range("J17").Value = [double]"0.0005804273921477662"
$ws.Range("M17").Value = [double]"79.39367866666667"
$ws.Range("N17").Value = [double]"238.181036"
$ws.Range("O17").Value = [double]"0.1295223675414243"
$ws.Range("P17").Value = [double]"0.1295223675414243"
$ws.Range("Q17").Value = [double]"225.3130938136236"
$ws.Range("R17").Value = [double]"2027.817844322612"
$ws.Range("S17").Value = [double]"7.517833001687341E-05"
$ws.Range("T17").Value = [double]"7.517833001687341E-05"
